$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

$changes = @(
    @{Addr='D2'; Val='67.313.89'},
    @{Addr='E2'; Val='  +2.31%  '},
    @{Addr='D3'; Val='3.368.17'},
    @{Addr='E3'; Val='  +1.78%  '},
    @{Addr='E4'; Val='  -0.25%  '},
    @{Addr='D5'; Val='589.52'},
    @{Addr='E5'; Val='  +6.43%  '},
    @{Addr='D6'; Val='188.68'},
    @{Addr='E6'; Val='  +0.79%  '},
    @{Addr='D7'; Val='0.999'},
    @{Addr='E7'; Val='  +0.03%  '},
    @{Addr='D8'; Val='0.597'},
    @{Addr='E8'; Val='  +3.21%  '},
    @{Addr='D9'; Val='0.185'},
    @{Addr='E9'; Val='  +3.47%  '},
    @{Addr='E10'; Val='  +1.27%  '},
    @{Addr='D11'; Val='47.78'},
    @{Addr='E11'; Val='  +2.38%  '},
    @{Addr='E12'; Val='  +3.25%  '},
    @{Addr='D13'; Val='659.55'},
    @{Addr='E13'; Val='  +10.44%  '},
    @{Addr='D14'; Val='3.907.10'},
    @{Addr='E14'; Val='  +1.81%  '},
    @{Addr='E15'; Val='  +0.40%  '},
    @{Addr='D16'; Val='67.354.07'},
    @{Addr='E16'; Val='  +2.34%  '},
    @{Addr='B17'; Val='Chainlink'},
    @{Addr='C17'; Val='https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'},
    @{Addr='D17'; Val='18.09'},
    @{Addr='E17'; Val='  +1.22%  '},
    @{Addr='B18'; Val='TRON'},
    @{Addr='C18'; Val='https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'},
    @{Addr='D18'; Val='0.119'},
    @{Addr='E18'; Val='  +1.06%  '},
    @{Addr='D19'; Val='3.372.21'},
    @{Addr='E19'; Val='  +1.99%  '},
    @{Addr='D20'; Val='11.21'},
    @{Addr='E20'; Val='  +2.15%  '},
    @{Addr='E21'; Val='  +1.65%  '},
    @{Addr='D22'; Val='18.09'},
    @{Addr='E22'; Val='  -1.99%  '},
    @{Addr='D23'; Val='5.10'},
    @{Addr='E23'; Val='  +1.31%  '},
    @{Addr='D24'; Val='100.97'},
    @{Addr='E24'; Val='  +0.85%  '},
    @{Addr='E25'; Val='  +2.74%  '},
    @{Addr='E26'; Val='  +4.47%  '},
    @{Addr='E27'; Val='  +3.75%  '},
    @{Addr='D28'; Val='32.30'},
    @{Addr='E28'; Val='  +6.53%  '},
    @{Addr='D29'; Val='8.76'},
    @{Addr='E29'; Val='  +1.22%  '},
    @{Addr='D30'; Val='6.85'},
    @{Addr='E30'; Val='  +2.73%  '},
    @{Addr='D31'; Val='621.26'},
    @{Addr='E31'; Val='  +9.24%  '},
    @{Addr='D32'; Val='3.90'},
    @{Addr='E32'; Val='  +1.63%  '},
    @{Addr='E33'; Val='  +2.27%  '},
    @{Addr='E34'; Val='  +3.46%  '},
    @{Addr='D35'; Val='3.891.62'},
    @{Addr='E35'; Val='  +5.01%  '},
    @{Addr='E36'; Val='  +0.27%  '},
    @{Addr='D37'; Val='55.70'},
    @{Addr='E37'; Val='  -1.71%  '},
    @{Addr='E38'; Val='  +7.87%  '},
    @{Addr='E39'; Val='  +2.86%  '},
    @{Addr='D40'; Val='33.81'},
    @{Addr='E40'; Val='  +0.95%  '},
    @{Addr='D41'; Val='3.28'},
    @{Addr='E41'; Val='  +2.13%  '},
    @{Addr='D42'; Val='0.0₃0708'},
    @{Addr='E42'; Val='  +1.80%  '},
    @{Addr='E43'; Val='  +3.35%  '},
    @{Addr='E44'; Val='  -0.89%  '},
    @{Addr='E45'; Val='  +2.30%  '},
    @{Addr='E46'; Val='  +1.19%  '},
    @{Addr='D47'; Val='2.61'},
    @{Addr='E47'; Val='  +2.70%  '},
    @{Addr='D49'; Val='2.88'},
    @{Addr='E49'; Val='  -18.38%  '},
    @{Addr='D50'; Val='1.36'},
    @{Addr='E50'; Val='  +10.56%  '},
    @{Addr='D51'; Val='130.48'},
    @{Addr='E51'; Val='  +5.22%  '}
)

foreach ($item in $changes) {
    Set-TextCell $ws $item.Addr $item.Val
}
